{"js": "// Adjust \"Figure Caption\" font to match the standard text font: the\n// \"Image Caption\" paragraph style (and its linked \"Image Caption Zchn\"\n// character style) switch from Arial to Times New Roman, per the\n// Springer submission guidelines referenced in the commit message.\n\nconst styles = context.document.getStyles();\n\nconst imageCaption = styles.getByNameOrNullObject(\"Image Caption\");\nconst imageCaptionChar = styles.getByNameOrNullObject(\"Image Caption Zchn\");\n\nimageCaption.load(\"isNullObject\");\nimageCaptionChar.load(\"isNullObject\");\nawait context.sync();\n\nif (!imageCaption.isNullObject) {\n  imageCaption.font.name = \"Times New Roman\";\n}\n\nif (!imageCaptionChar.isNullObject) {\n  imageCaptionChar.font.name = \"Times New Roman\";\n}\n\nawait context.sync();\n", "ps1": "# Adjust \"Figure Caption\" font to match the standard text font: the\n# \"Image Caption\" paragraph style (and its linked \"Image Caption Zchn\"\n# character style) switch from Arial to Times New Roman, per the\n# Springer submission guidelines referenced in the commit message.\n\n$d = $word.ActiveDocument\n\n$d.Styles(\"Image Caption\").Font.Name = \"Times New Roman\"\n$d.Styles(\"Image Caption Zchn\").Font.Name = \"Times New Roman\"\n"}
